# Update cryptocurrency price and 1h volume/change figures in the "cryptos" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to remain plain text (matches original inlineStr cells)
    # instead of letting Excel auto-convert numeric-looking strings to numbers,
    # then restore the default (unstyled) cell style so no extra formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.294.71"
Set-TextValue $ws.Range("E2") "  -2.74%  "

Set-TextValue $ws.Range("D3") "3.378.96"
Set-TextValue $ws.Range("E3") "  -2.10%  "

Set-TextValue $ws.Range("E4") "  +0.02%  "

Set-TextValue $ws.Range("D5") "573.10"
Set-TextValue $ws.Range("E5") "  -0.80%  "

Set-TextValue $ws.Range("D6") "151.80"
Set-TextValue $ws.Range("E6") "  +1.77%  "

Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("E8") "  +1.12%  "

Set-TextValue $ws.Range("D9") "8.01"
Set-TextValue $ws.Range("E9") "  +2.50%  "

Set-TextValue $ws.Range("E10") "  -0.42%  "

Set-TextValue $ws.Range("E11") "  +3.12%  "

Set-TextValue $ws.Range("D12") "3.961.62"
Set-TextValue $ws.Range("E12") "  -1.92%  "

Set-TextValue $ws.Range("E13") "  +1.02%  "

Set-TextValue $ws.Range("D14") "28.44"
Set-TextValue $ws.Range("E14") "  -0.64%  "

Set-TextValue $ws.Range("E15") "  -0.21%  "

Set-TextValue $ws.Range("D16") "3.399.41"
Set-TextValue $ws.Range("E16") "  -1.41%  "

Set-TextValue $ws.Range("D17") "61.367.41"
Set-TextValue $ws.Range("E17") "  -2.63%  "

Set-TextValue $ws.Range("E18") "  -0.67%  "

Set-TextValue $ws.Range("D19") "14.28"
Set-TextValue $ws.Range("E19") "  -0.88%  "

Set-TextValue $ws.Range("D20") "8.89"
Set-TextValue $ws.Range("E20") "  -2.62%  "

Set-TextValue $ws.Range("D21") "375.95"
Set-TextValue $ws.Range("E21") "  -2.37%  "

Set-TextValue $ws.Range("D22") "0.567"
Set-TextValue $ws.Range("E22") "  +1.26%  "

Set-TextValue $ws.Range("D23") "75.38"
Set-TextValue $ws.Range("E23") "  +1.31%  "

Set-TextValue $ws.Range("E24") "  -0.01%  "

Set-TextValue $ws.Range("D25") "3.543.07"
Set-TextValue $ws.Range("E25") "  -1.08%  "

Set-TextValue $ws.Range("E26") "  -4.82%  "

Set-TextValue $ws.Range("E27") "  -3.73%  "

Set-TextValue $ws.Range("E28") "  -2.76%  "

Set-TextValue $ws.Range("E29") "  +0.34%  "

Set-TextValue $ws.Range("E30") "  +0.57%  "

Set-TextValue $ws.Range("E31") "  -0.01%  "

Set-TextValue $ws.Range("D32") "7.72"
Set-TextValue $ws.Range("E32") "  -3.62%  "

Set-TextValue $ws.Range("D33") "23.09"
Set-TextValue $ws.Range("E33") "  -0.82%  "

Set-TextValue $ws.Range("E34") "  -2.15%  "

Set-TextValue $ws.Range("E35") "  +0.80%  "

Set-TextValue $ws.Range("D36") "169.72"
Set-TextValue $ws.Range("E36") "  -0.03%  "

Set-TextValue $ws.Range("E37") "  -3.85%  "

Set-TextValue $ws.Range("D38") "6.81"
Set-TextValue $ws.Range("E38") "  -3.23%  "

Set-TextValue $ws.Range("D39") "30.43"
Set-TextValue $ws.Range("E39") "  -4.31%  "

Set-TextValue $ws.Range("D40") "3.413.74"
Set-TextValue $ws.Range("E40") "  -2.07%  "

Set-TextValue $ws.Range("E41") "  +0.17%  "

Set-TextValue $ws.Range("D42") "42.50"
Set-TextValue $ws.Range("E42") "  +0.10%  "

Set-TextValue $ws.Range("E43") "  -2.71%  "

Set-TextValue $ws.Range("E44") "  -0.12%  "

Set-TextValue $ws.Range("E45") "  -6.14%  "

Set-TextValue $ws.Range("E46") "  -4.35%  "

Set-TextValue $ws.Range("D47") "2.522.10"
Set-TextValue $ws.Range("E47") "  -2.20%  "

Set-TextValue $ws.Range("D48") "23.01"
Set-TextValue $ws.Range("E48") "  +2.06%  "

Set-TextValue $ws.Range("E49") "  -1.39%  "

Set-TextValue $ws.Range("E50") "  +0.05%  "

Set-TextValue $ws.Range("E51") "  -1.83%  "
